# Auto commit at 2025-11-22  8:06:23.19
# Append the day's (2025-11-21 / Excel serial 45982) readings for both
# charging stations to the daily data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Last existing data row is row 41 ("高岭站" for 2025-11-20). Grab its
# formatting (date format on A, "0.00" format on C:E, "0" format on F)
# so the two freshly appended rows (42 and 43) keep the same look instead
# of picking up the workbook's default "General" style.
$ws.Range("A41:F41").Copy()
$ws.Range("A42:F43").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 42: 四方坪站 (Sifangping station)
$ws.Cells.Item(42, 1).Value = 45982
$ws.Cells.Item(42, 2).Value = "四方坪站"
$ws.Cells.Item(42, 3).Value = 9400.8700000000008
$ws.Cells.Item(42, 4).Value = 8398.91
$ws.Cells.Item(42, 5).Value = 3134.01
$ws.Cells.Item(42, 6).Value = 407

# Row 43: 高岭站 (Gaoling station)
$ws.Cells.Item(43, 1).Value = 45982
$ws.Cells.Item(43, 2).Value = "高岭站"
$ws.Cells.Item(43, 3).Value = 4978.18
$ws.Cells.Item(43, 4).Value = 4302.78
$ws.Cells.Item(43, 5).Value = 1290.75
$ws.Cells.Item(43, 6).Value = 173

# Reflect the scrolled/selected state captured when the edit was made:
# the view had scrolled down so row 34 is at the top, with J45 selected.
$win = $excel.ActiveWindow
$win.ScrollRow = 34
$win.ScrollColumn = 1
$ws.Range("J45").Select()
